# Apply the cryptocurrency price/volume updates described by the commit diff.
# Column D ("Price") and column E ("Volume(1h)") cells hold formatted, text-typed
# values (e.g. "63.690.01", "  -0.77%  "). Values that look like a plain number
# (e.g. "591.46") would be auto-converted to a numeric cell by Excel when assigned
# through .Value, so those are entered with a leading apostrophe (Excel's classic
# "force text" quote-prefix) to keep them as text, matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "63.690.01"
    "E2" = "  -0.77%  "
    "D3" = "2.615.13"
    "E3" = "  -0.30%  "
    "E4" = "  -0.01%  "
    "D5" = "'591.46"
    "E5" = "  -1.75%  "
    "D6" = "'150.23"
    "E6" = "  -0.12%  "
    "E7" = "  +0.02%  "
    "E8" = "  -0.65%  "
    "E9" = "  +0.34%  "
    "D10" = "'5.77"
    "E10" = "  +1.27%  "
    "D11" = "'0.389"
    "E11" = "  +0.83%  "
    "E12" = "  +0.43%  "
    "D13" = "'27.76"
    "E13" = "  +0.37%  "
    "D14" = "3.084.32"
    "E14" = "  -0.32%  "
    "D15" = "63.481.94"
    "E15" = "  -0.84%  "
    "E16" = "  +5.01%  "
    "D17" = "2.637.94"
    "E17" = "  +1.42%  "
    "E18" = "  +0.14%  "
    "D19" = "'4.78"
    "E19" = "  +2.69%  "
    "D20" = "'346.42"
    "E20" = "  -1.35%  "
    "D21" = "'6.91"
    "E21" = "  -0.49%  "
    "E22" = "  +0.18%  "
    "D23" = "'67.12"
    "E23" = "  +0.91%  "
    "E24" = "  -2.98%  "
    "D25" = "'9.25"
    "E25" = "  +0.12%  "
    "D26" = "'1.67"
    "E26" = "  -1.25%  "
    "D27" = "'8.57"
    "E27" = "  +4.41%  "
    "D28" = "'547.93"
    "E28" = "  +1.56%  "
    "E29" = "  -2.02%  "
    "E30" = "  -0.20%  "
    "E31" = "  -1.02%  "
    "D32" = "0.0₃0868"
    "E32" = "  +1.60%  "
    "D33" = "'1.79"
    "E33" = "  +1.88%  "
    "D34" = "'5.37"
    "E34" = "  +1.23%  "
    "D35" = "'6.16"
    "E35" = "  +0.42%  "
    "D36" = "'165.70"
    "E36" = "  -1.28%  "
    "E37" = "  +0.97%  "
    "E38" = "  -0.03%  "
    "E39" = "  -1.06%  "
    "D40" = "'19.57"
    "E40" = "  +0.43%  "
    "D42" = "'165.82"
    "E42" = "  -2.20%  "
    "D43" = "'4.10"
    "E43" = "  +4.22%  "
    "D44" = "'23.29"
    "E44" = "  +7.69%  "
    "D45" = "'0.0583"
    "E45" = "  -1.62%  "
    "E46" = "  +7.79%  "
    "E48" = "  +1.79%  "
    "D49" = "'0.0961"
    "E49" = "  -0.76%  "
    "D50" = "'19.25"
    "E50" = "  -0.30%  "
    "D51" = "0.0₆0232"
    "E51" = "  +17.50%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
